$d = $word.ActiveDocument

# Remove the " 3" text (third run) entirely.
$d.Content.Find.Execute(" 3", $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 2)
